# Adds a new "2023" data column (S) plus the trailing blank formatting
# column (U) that Excel keeps after the last populated column, mirroring
# the existing 2022 column (R) / trailing blank column (T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new year header (2023) in S3, copying R3's look ---
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S3").Value = 2023

# --- Rows 4-14: new data values in column S, copying column R's style ---
$values = @{
    4  = 1383.1
    5  = 228.6
    6  = 7469
    7  = 5226
    8  = 724.1
    9  = 56.3
    10 = 46.2
    11 = 172.2
    12 = 13
    13 = 1190.7
    14 = 1145.9000000000001
}

foreach ($row in 4..14) {
    $ws.Range("R$row").Copy()
    $ws.Range("S$row").PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("S$row").Value = $values[$row]
}

# --- New trailing blank column U, copying column T's formatting ---
# Row 2 (header band), row 3 (year band), rows 4-14 (data band) and
# row 15 (footer band) all already carry a blank styled cell in T;
# replicate that same style one column over into U.
foreach ($row in 2..15) {
    $ws.Range("T$row").Copy()
    $ws.Range("U$row").PasteSpecial(-4122)   # xlPasteFormats
}

# --- Row 1: extend the title bar's shaded style into the new S column ---
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)   # xlPasteFormats

# --- Merge the title cell across the now-wider A1:S1 range ---
$ws.Range("A1:R1").UnMerge()
$ws.Range("A1:S1").Merge()

# --- Update the selection / view to match the new layout ---
$ws.Range("I8").Select()
$ws.Range("S3:S14").Select()

$excel.CutCopyMode = $false
